$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Updated forecast error values (bugfix for selection issue in naive_dict)
$ws.Range("B2").Value = 0.1420281398799347
$ws.Range("C2").Value = 0.96300297170131
$ws.Range("D2").Value = 4.425279606777181
$ws.Range("E2").Value = 2.103634855857162
$ws.Range("F2").Value = 2.119719265974769
$ws.Range("G2").Value = 51

$ws.Range("B3").Value = 0.1234272948742712
$ws.Range("C3").Value = 1.013190509142516
$ws.Range("D3").Value = 4.456324062957081
$ws.Range("E3").Value = 2.111000725475261
$ws.Range("F3").Value = 2.128784688702496
$ws.Range("G3").Value = 50

$ws.Range("B4").Value = 0.1607615939753095
$ws.Range("C4").Value = 0.931089636352725
$ws.Range("D4").Value = 4.199434463388394
$ws.Range("E4").Value = 2.049252171741778
$ws.Range("F4").Value = 2.064107547723181
$ws.Range("G4").Value = 49

$ws.Range("B5").Value = 0.1705283119469347
$ws.Range("C5").Value = 1.039568901588029
$ws.Range("D5").Value = 4.679721022412824
$ws.Range("E5").Value = 2.163266285599816
$ws.Range("F5").Value = 2.17935562225223
$ws.Range("G5").Value = 48

$ws.Range("B6").Value = 0.1525391686516503
$ws.Range("C6").Value = 0.9855784814025781
$ws.Range("D6").Value = 4.478113517829056
$ws.Range("E6").Value = 2.116155362403492
$ws.Range("F6").Value = 2.133468955323994
$ws.Range("G6").Value = 47

$ws.Range("B7").Value = 0.164531530408786
$ws.Range("C7").Value = 1.008196845425907
$ws.Range("D7").Value = 4.671313786093526
$ws.Range("E7").Value = 2.16132223097194
$ws.Range("F7").Value = 2.178864053856015
$ws.Range("G7").Value = 46

$ws.Range("B8").Value = 0.1169720003859123
$ws.Range("C8").Value = 0.9369220677651683
$ws.Range("D8").Value = 4.351949521053745
$ws.Range("E8").Value = 2.086132671009624
$ws.Range("F8").Value = 2.106386491385813
$ws.Range("G8").Value = 45

$wb.Save()
